$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 NumOrden first: 1078703 (text)
$ws.Range("D4").Value = "'1078703"

# Row 4: QA / 1220194200684 (trailing space) / Importe cleared
$ws.Range("B4").Value = "'1220194200684 "
$ws.Range("C4").ClearContents()

# Row 3: QA / 1120194100442 (trailing space) / Importe cleared
$ws.Range("B3").Value = "'1120194100442 "
$ws.Range("C3").ClearContents()

# Row 3 NumOrden: 1078702 (text)
$ws.Range("D3").Value = "'1078702"

# Row 2: QA / 0420194406830 (no trailing space) / Importe cleared / NumOrden empty
$ws.Range("B2").Value = "'0420194406830"
$ws.Range("C2").ClearContents()

# Row 5: PREPROD / 0420172008629 (trailing space) / Importe cleared / NumOrden 1703284 (number, quote-prefixed style)
$ws.Range("B5").Value = "'0420172008629 "
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 1703284
$ws.Range("B5").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Row 6: PREPROD / 1120170200942 (trailing space) / Importe cleared / NumOrden 1703286 (number, quote-prefixed style)
$ws.Range("B6").Value = "'1120170200942 "
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 1703286
$ws.Range("B6").Copy()
$ws.Range("D6").PasteSpecial(-4122)

# Row 7: PREPROD / 1220170301442 (trailing space) / Importe cleared / NumOrden 1703285 (number, quote-prefixed style)
$ws.Range("B7").Value = "'1220170301442 "
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 1703285
$ws.Range("B7").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# D1 header text stays "NumOrden" (unaffected in content, only shared string index shifts automatically)

# Update the active selection to reflect the last worked cell
$ws.Range("E6").Select()

$excel.CutCopyMode = 0
